$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("D2").Value = 8.4700000000000006
$ws.Range("D3").Value = 8.76
$ws.Range("D4").Value = 6.8
$ws.Range("D5").Value = 6.53
$ws.Range("D6").Value = 4.4000000000000004
$ws.Range("D7").Value = 4.4000000000000004
$ws.Range("D8").Value = 4.21
$ws.Range("D9").Value = 4.26
$ws.Range("D10").Value = 4.16
$ws.Range("D11").Value = 4.09
$ws.Range("D12").Value = 4.3100000000000005
$ws.Range("D13").Value = 3.56
$ws.Range("D14").Value = ""
$ws.Range("D15").Value = 4.16
$ws.Range("D16").Value = 4.5
$ws.Range("D17").Value = 7.62
$ws.Range("D18").Value = 7.24
$ws.Range("D19").Value = 8.81

$ws.Range("D14").Select()
